$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each price (column D) and 1h-volume % (column E) cell holds a text-
# formatted value (t="inlineStr" in the OOXML). We prefix with a single
# quote to force Excel to store the new value as text rather than re-
# interpreting it as a number/percentage, then reset the style back to
# "Normal" so no stray text/quote-prefix number format sticks to the cell.

$ws.Range("D2").Value = "'274.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.46%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.73%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.852"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.11%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06324"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.904"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.09%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'1.82%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.273"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'34.97%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8720"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.73%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'0.74%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.04989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.46%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.17%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02954"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-6.48%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09029"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.54%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006312"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.63%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006028"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.41%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.13%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.284"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.10%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1334"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.73%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.906"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.58%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04360"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.94%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.20%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004253"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.53%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-0.23%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04052"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.59%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006713"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.54%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1165"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.05%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.29%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01070"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-10.77%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005296"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.98%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.486"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-37.29%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.02000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-33.03%"
$ws.Range("E47").Style = "Normal"
